$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update payout values (single-zero wheel: 37 sectors instead of 38)
$ws.Range("B4").Value = 38
$ws.Range("B9").Value = 4

# Update the probability formula denominator from 38 to 37.
# D4 has its own formula; D5:D11 is a shared formula group.
$ws.Range("D4").Formula = "=C4/37"
$ws.Range("D5:D11").Formula = "=C5/37"

# Restore the active cell selection
$ws.Range("B10").Select()
